$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Internal Assignment" column (O) values from the export.
# Keep cell formatting/styles intact, only clear the contents so the
# now-unused "Internal Assignment" shared string gets dropped.
$ws.Range("O4").ClearContents()
$ws.Range("O5:O7").ClearContents()
$ws.Range("O12:O15").ClearContents()
$ws.Range("O20:O25").ClearContents()

# Update the active selection to match the edited column.
$ws.Range("O4:O25").Select()
